# New test run: update the start/end date & time columns on the single
# data row (row 2) of the Katalon test-data sheet to the new run's
# timestamps, leaving everything else (headers, other columns, styles)
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "16/10/2019"      # Fecha Inicio
$ws.Range("E2").Value = "09:39:18.263"    # Hora Inicio
$ws.Range("F2").Value = "16/10/2019"      # Fecha Final
$ws.Range("G2").Value = "09:40:29.013"    # Hora Final
